$wb = $excel.ActiveWorkbook

# --- Sheet: average_mae ---
$ws1 = $wb.Worksheets.Item("average_mae")
$ws1.Cells.Item(5, 1).Value = "ibes_1|ni-sector_code|cnn_rnn｜all"
$ws1.Cells.Item(6, 1).Value = "ibes_1|fwdepsqcut|dense2｜all x 0 -fix space"
$ws1.Cells.Item(7, 1).Value = "ibes_6|fwdepsqcut|ibes_sector_only ws -indi space"
$ws1.Cells.Item(8, 1).Value = "ibes_1|fwdepsqcut|ibes_entire_only ws -smaller space"
$ws1.Cells.Item(9, 1).Value = "ibes_1|fwdepsqcut-industry_code|ibes_entire_only ws -smaller space"
$ws1.Cells.Item(10, 1).Value = "ibes_1|fwdepsqcut-sector_code|ibes_entire_only ws -smaller space"
$ws1.Cells.Item(11, 1).Value = "ibes_2|ni|ibes_new industry_all x -indi space"
$ws1.Cells.Item(12, 1).Value = "ibes_2|fwdepsqcut|dense2｜new industry model -fix space"
$ws1.Cells.Item(13, 1).Value = "ibes_1|fwdepsqcut|dense2｜new with indi code -fix space"
$ws1.Cells.Item(14, 1).Value = "ibes_1|fwdepsqcut-industry_code|dense2｜new with indi code -fix space"
$ws1.Cells.Item(15, 1).Value = "ibes_1|fwdepsqcut-sector_code|dense2｜new with indi code -fix space"
$ws1.Cells.Item(16, 1).Value = "ibes_1|fwdepsqcut|rnn_eps｜all"
$ws1.Cells.Item(4, 2).Value = 0.00870510303297917
$ws1.Cells.Item(4, 3).Value = 0.009208389006560668
$ws1.Cells.Item(4, 4).Value = 12746
$ws1.Cells.Item(5, 2).Value = 0.00873659052644278
$ws1.Cells.Item(5, 3).Value = 0.009050386758606433
$ws1.Cells.Item(5, 4).Value = 12054
$ws1.Cells.Item(6, 2).Value = 0.008861431855592515
$ws1.Cells.Item(6, 3).Value = 0.009338113103833704
$ws1.Cells.Item(6, 4).Value = 15176
$ws1.Cells.Item(7, 2).Value = 0.008875614826370663
$ws1.Cells.Item(7, 3).Value = 0.008557495853986283
$ws1.Cells.Item(7, 4).Value = 15176
$ws1.Cells.Item(8, 2).Value = 0.008861431855592515
$ws1.Cells.Item(8, 3).Value = 0.00865559650904604
$ws1.Cells.Item(8, 4).Value = 15176
$ws1.Cells.Item(9, 2).Value = 0.008861431855592515
$ws1.Cells.Item(9, 3).Value = 0.008658869869704775
$ws1.Cells.Item(9, 4).Value = 15176
$ws1.Cells.Item(10, 2).Value = 0.008861431855592515
$ws1.Cells.Item(10, 3).Value = 0.008630337442592475
$ws1.Cells.Item(10, 4).Value = 15176
$ws1.Cells.Item(11, 2).Value = 0.008828774570773384
$ws1.Cells.Item(11, 3).Value = 0.007811737484649521
$ws1.Cells.Item(11, 4).Value = 15176
$ws1.Cells.Item(12, 2).Value = 0.008828774570773384
$ws1.Cells.Item(12, 3).Value = 0.009833523580700576
$ws1.Cells.Item(12, 4).Value = 15176
$ws1.Cells.Item(13, 2).Value = 0.008861431855592515
$ws1.Cells.Item(13, 3).Value = 0.009476279475976703
$ws1.Cells.Item(13, 4).Value = 15176
$ws1.Cells.Item(14, 2).Value = 0.008861431855592515
$ws1.Cells.Item(14, 3).Value = 0.009848581092865606
$ws1.Cells.Item(14, 4).Value = 15176
$ws1.Cells.Item(15, 2).Value = 0.008861431855592515
$ws1.Cells.Item(15, 3).Value = 0.01057408267989892
$ws1.Cells.Item(15, 4).Value = 15176
$ws1.Cells.Item(16, 2).Value = 0.0087378846049827
$ws1.Cells.Item(16, 3).Value = 0.009921528473829913
$ws1.Cells.Item(16, 4).Value = 6584

# --- Sheet: average_mse ---
$ws2 = $wb.Worksheets.Item("average_mse")
$ws2.Cells.Item(5, 1).Value = "ibes_1|ni-sector_code|cnn_rnn｜all"
$ws2.Cells.Item(6, 1).Value = "ibes_1|fwdepsqcut|dense2｜all x 0 -fix space"
$ws2.Cells.Item(7, 1).Value = "ibes_6|fwdepsqcut|ibes_sector_only ws -indi space"
$ws2.Cells.Item(8, 1).Value = "ibes_1|fwdepsqcut|ibes_entire_only ws -smaller space"
$ws2.Cells.Item(9, 1).Value = "ibes_1|fwdepsqcut-industry_code|ibes_entire_only ws -smaller space"
$ws2.Cells.Item(10, 1).Value = "ibes_1|fwdepsqcut-sector_code|ibes_entire_only ws -smaller space"
$ws2.Cells.Item(11, 1).Value = "ibes_2|ni|ibes_new industry_all x -indi space"
$ws2.Cells.Item(12, 1).Value = "ibes_2|fwdepsqcut|dense2｜new industry model -fix space"
$ws2.Cells.Item(13, 1).Value = "ibes_1|fwdepsqcut|dense2｜new with indi code -fix space"
$ws2.Cells.Item(14, 1).Value = "ibes_1|fwdepsqcut-industry_code|dense2｜new with indi code -fix space"
$ws2.Cells.Item(15, 1).Value = "ibes_1|fwdepsqcut-sector_code|dense2｜new with indi code -fix space"
$ws2.Cells.Item(16, 1).Value = "ibes_1|fwdepsqcut|rnn_eps｜all"
$ws2.Cells.Item(4, 2).Value = 0.0001918729794254194
$ws2.Cells.Item(4, 3).Value = 0.0001801854361867948
$ws2.Cells.Item(5, 2).Value = 0.0001944778043602815
$ws2.Cells.Item(5, 3).Value = 0.0001765025584832404
$ws2.Cells.Item(6, 2).Value = 0.0001975833946462478
$ws2.Cells.Item(6, 3).Value = 0.0001881828907176001
$ws2.Cells.Item(7, 2).Value = 0.0002043910143630977
$ws2.Cells.Item(7, 3).Value = 0.0001720965334955223
$ws2.Cells.Item(8, 2).Value = 0.0001975833946462478
$ws2.Cells.Item(8, 3).Value = 0.0001658971970851478
$ws2.Cells.Item(9, 2).Value = 0.0001975833946462478
$ws2.Cells.Item(9, 3).Value = 0.0001657197735352102
$ws2.Cells.Item(10, 2).Value = 0.0001975833946462478
$ws2.Cells.Item(10, 3).Value = 0.0001649945689266668
$ws2.Cells.Item(11, 2).Value = 0.0002002415943724616
$ws2.Cells.Item(11, 3).Value = 0.0001454923001376933
$ws2.Cells.Item(12, 2).Value = 0.0002002415943724616
$ws2.Cells.Item(12, 3).Value = 0.0002242374076787091
$ws2.Cells.Item(13, 2).Value = 0.0001975833946462478
$ws2.Cells.Item(13, 3).Value = 0.0001902112594253567
$ws2.Cells.Item(14, 2).Value = 0.0001975833946462478
$ws2.Cells.Item(14, 3).Value = 0.0002044223000088748
$ws2.Cells.Item(15, 2).Value = 0.0001975833946462478
$ws2.Cells.Item(15, 3).Value = 0.0002335431519097555
$ws2.Cells.Item(16, 2).Value = 0.000196008906973461
$ws2.Cells.Item(16, 3).Value = 0.0002236856401977318

# --- Sheet: average_r2 ---
$ws3 = $wb.Worksheets.Item("average_r2")
$ws3.Cells.Item(5, 1).Value = "ibes_1|ni-sector_code|cnn_rnn｜all"
$ws3.Cells.Item(6, 1).Value = "ibes_1|fwdepsqcut|dense2｜all x 0 -fix space"
$ws3.Cells.Item(7, 1).Value = "ibes_6|fwdepsqcut|ibes_sector_only ws -indi space"
$ws3.Cells.Item(8, 1).Value = "ibes_1|fwdepsqcut|ibes_entire_only ws -smaller space"
$ws3.Cells.Item(9, 1).Value = "ibes_1|fwdepsqcut-industry_code|ibes_entire_only ws -smaller space"
$ws3.Cells.Item(10, 1).Value = "ibes_1|fwdepsqcut-sector_code|ibes_entire_only ws -smaller space"
$ws3.Cells.Item(11, 1).Value = "ibes_2|ni|ibes_new industry_all x -indi space"
$ws3.Cells.Item(12, 1).Value = "ibes_2|fwdepsqcut|dense2｜new industry model -fix space"
$ws3.Cells.Item(13, 1).Value = "ibes_1|fwdepsqcut|dense2｜new with indi code -fix space"
$ws3.Cells.Item(14, 1).Value = "ibes_1|fwdepsqcut-industry_code|dense2｜new with indi code -fix space"
$ws3.Cells.Item(15, 1).Value = "ibes_1|fwdepsqcut-sector_code|dense2｜new with indi code -fix space"
$ws3.Cells.Item(16, 1).Value = "ibes_1|fwdepsqcut|rnn_eps｜all"
$ws3.Cells.Item(4, 2).Value = 0.1261272704213072
$ws3.Cells.Item(4, 3).Value = 0.2601499020947659
$ws3.Cells.Item(4, 4).Value = 0.1793574091442834
$ws3.Cells.Item(5, 2).Value = 0.1185037215870294
$ws3.Cells.Item(5, 3).Value = 0.260106022161864
$ws3.Cells.Item(5, 4).Value = 0.1999788924749928
$ws3.Cells.Item(6, 2).Value = 0.1445666443086192
$ws3.Cells.Item(6, 3).Value = 0.2636440167903586
$ws3.Cells.Item(6, 4).Value = 0.1852659380689607
$ws3.Cells.Item(7, 2).Value = 0.1726321527164353
$ws3.Cells.Item(7, 3).Value = 0.2636440167903589
$ws3.Cells.Item(7, 4).Value = 0.3033591085848544
$ws3.Cells.Item(8, 2).Value = 0.1445666443086192
$ws3.Cells.Item(8, 3).Value = 0.2636440167903586
$ws3.Cells.Item(8, 4).Value = 0.2817514029636731
$ws3.Cells.Item(9, 2).Value = 0.1445666443086192
$ws3.Cells.Item(9, 3).Value = 0.2636440167903586
$ws3.Cells.Item(9, 4).Value = 0.2825195546748713
$ws3.Cells.Item(10, 2).Value = 0.1445666443086192
$ws3.Cells.Item(10, 3).Value = 0.2636440167903586
$ws3.Cells.Item(10, 4).Value = 0.2856593135243424
$ws3.Cells.Item(11, 2).Value = 0.1712998665722781
$ws3.Cells.Item(11, 3).Value = 0.2636440167903586
$ws3.Cells.Item(11, 4).Value = 0.3978799014527122
$ws3.Cells.Item(12, 2).Value = 0.1712998665722781
$ws3.Cells.Item(12, 3).Value = 0.2636440167903586
$ws3.Cells.Item(12, 4).Value = 0.07199315783919602
$ws3.Cells.Item(13, 2).Value = 0.1445666443086192
$ws3.Cells.Item(13, 3).Value = 0.2636440167903586
$ws3.Cells.Item(13, 4).Value = 0.1764841563136554
$ws3.Cells.Item(14, 2).Value = 0.1445666443086192
$ws3.Cells.Item(14, 3).Value = 0.2636440167903586
$ws3.Cells.Item(14, 4).Value = 0.114957740310983
$ws3.Cells.Item(15, 2).Value = 0.1445666443086192
$ws3.Cells.Item(15, 3).Value = 0.2636440167903586
$ws3.Cells.Item(15, 4).Value = -0.01112040561197025
$ws3.Cells.Item(16, 2).Value = 0.1037236563029756
$ws3.Cells.Item(16, 3).Value = 0.3312793915445512
$ws3.Cells.Item(16, 4).Value = -0.02283182345941115
